$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Datos actualizados a 13 de Junio de 2020 a las 17:17
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 13 de Junio de 2020 a las 17:17'

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = 'Estados Unidos'
$ws.Cells.Item(4, 2).Value = 2122343
$ws.Cells.Item(4, 3).Value = 5421
$ws.Cells.Item(4, 4).Value = 842308
$ws.Cells.Item(4, 5).Value = 1163106
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 104
$ws.Cells.Item(4, 8).Value = 116929

# Row 7: India
$ws.Cells.Item(7, 1).Value = 'India'
$ws.Cells.Item(7, 2).Value = 310760
$ws.Cells.Item(7, 3).Value = 1157
$ws.Cells.Item(7, 4).Value = 155290
$ws.Cells.Item(7, 5).Value = 146575
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 5
$ws.Cells.Item(7, 8).Value = 8895

# Row 15: Chile
$ws.Cells.Item(15, 1).Value = 'Chile'
$ws.Cells.Item(15, 2).Value = 167355
$ws.Cells.Item(15, 3).Value = 6509
$ws.Cells.Item(15, 4).Value = 131358
$ws.Cells.Item(15, 5).Value = 32896
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 231
$ws.Cells.Item(15, 8).Value = 3101

# Row 19: Arabia Saudita
$ws.Cells.Item(19, 1).Value = 'Arabia Saudita'
$ws.Cells.Item(19, 2).Value = 123308
$ws.Cells.Item(19, 3).Value = 3366
$ws.Cells.Item(19, 4).Value = 82548
$ws.Cells.Item(19, 5).Value = 39828
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 39
$ws.Cells.Item(19, 8).Value = 932

# Row 33: Singapur
$ws.Cells.Item(33, 1).Value = 'Singapur'
$ws.Cells.Item(33, 2).Value = 40197
$ws.Cells.Item(33, 3).Value = 347
$ws.Cells.Item(33, 4).Value = 28808
$ws.Cells.Item(33, 5).Value = 11364
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 0
$ws.Cells.Item(33, 8).Value = 25

# Row 35: Portugal
$ws.Cells.Item(35, 1).Value = 'Portugal'
$ws.Cells.Item(35, 2).Value = 36463
$ws.Cells.Item(35, 3).Value = 283
$ws.Cells.Item(35, 4).Value = 22438
$ws.Cells.Item(35, 5).Value = 12513
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = 7
$ws.Cells.Item(35, 8).Value = 1512

# Row 40: Argentina
$ws.Cells.Item(40, 1).Value = 'Argentina'
$ws.Cells.Item(40, 2).Value = 28764
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = 9083
$ws.Cells.Item(40, 5).Value = 18879
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 17
$ws.Cells.Item(40, 8).Value = 802

# Row 44: Republica Dominicana
$ws.Cells.Item(44, 1).Value = 'Republica Dominicana'
$ws.Cells.Item(44, 2).Value = 22572
$ws.Cells.Item(44, 3).Value = 564
$ws.Cells.Item(44, 4).Value = 13084
$ws.Cells.Item(44, 5).Value = 8911
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 9
$ws.Cells.Item(44, 8).Value = 577

# Row 45: Oman
$ws.Cells.Item(45, 1).Value = 'Oman'
$ws.Cells.Item(45, 2).Value = 22077
$ws.Cells.Item(45, 3).Value = 1006
$ws.Cells.Item(45, 4).Value = 7530
$ws.Cells.Item(45, 5).Value = 14448
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(45, 7).Value = 3
$ws.Cells.Item(45, 8).Value = 99

# Row 48: Irak
$ws.Cells.Item(48, 1).Value = 'Irak'
$ws.Cells.Item(48, 2).Value = 18950
$ws.Cells.Item(48, 3).Value = 1180
$ws.Cells.Item(48, 4).Value = 7515
$ws.Cells.Item(48, 5).Value = 10886
$ws.Cells.Item(48, 6).Value = 0
$ws.Cells.Item(48, 7).Value = 53
$ws.Cells.Item(48, 8).Value = 549

# Row 49: Israel
$ws.Cells.Item(49, 1).Value = 'Israel'
$ws.Cells.Item(49, 2).Value = 18876
$ws.Cells.Item(49, 3).Value = 81
$ws.Cells.Item(49, 4).Value = 15319
$ws.Cells.Item(49, 5).Value = 3257
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 300

# Row 56: Kazajistan
$ws.Cells.Item(56, 1).Value = 'Kazajistan'
$ws.Cells.Item(56, 2).Value = 14238
$ws.Cells.Item(56, 3).Value = 366
$ws.Cells.Item(56, 4).Value = 8829
$ws.Cells.Item(56, 5).Value = 5337
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(56, 7).Value = 2
$ws.Cells.Item(56, 8).Value = 72

# Row 57: Serbia
$ws.Cells.Item(57, 1).Value = 'Serbia'
$ws.Cells.Item(57, 2).Value = 12251
$ws.Cells.Item(57, 3).Value = 76
$ws.Cells.Item(57, 4).Value = 11348
$ws.Cells.Item(57, 5).Value = 650
$ws.Cells.Item(57, 6).Value = 0
$ws.Cells.Item(57, 7).Value = 1
$ws.Cells.Item(57, 8).Value = 253

# Row 60: Moldavia
$ws.Cells.Item(60, 1).Value = 'Moldavia'
$ws.Cells.Item(60, 2).Value = 11459
$ws.Cells.Item(60, 3).Value = 366
$ws.Cells.Item(60, 4).Value = 6421
$ws.Cells.Item(60, 5).Value = 4640
$ws.Cells.Item(60, 6).Value = 0
$ws.Cells.Item(60, 7).Value = 13
$ws.Cells.Item(60, 8).Value = 398

# Row 61: Ghana
$ws.Cells.Item(61, 1).Value = 'Ghana'
$ws.Cells.Item(61, 2).Value = 11118
$ws.Cells.Item(61, 3).Value = 262
$ws.Cells.Item(61, 4).Value = 3979
$ws.Cells.Item(61, 5).Value = 7091
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 48

# Row 68: Noruega
$ws.Cells.Item(68, 1).Value = 'Noruega'
$ws.Cells.Item(68, 2).Value = 8625
$ws.Cells.Item(68, 3).Value = 5
$ws.Cells.Item(68, 4).Value = 8138
$ws.Cells.Item(68, 5).Value = 245
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 242

# Row 76: Tayikistan
$ws.Cells.Item(76, 1).Value = 'Tayikistan'
$ws.Cells.Item(76, 2).Value = 4971
$ws.Cells.Item(76, 3).Value = 69
$ws.Cells.Item(76, 4).Value = 3288
$ws.Cells.Item(76, 5).Value = 1633
$ws.Cells.Item(76, 6).Value = 0
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = 50

# Row 77: Uzbekistan
$ws.Cells.Item(77, 1).Value = 'Uzbekistan'
$ws.Cells.Item(77, 2).Value = 4937
$ws.Cells.Item(77, 3).Value = 68
$ws.Cells.Item(77, 4).Value = 3837
$ws.Cells.Item(77, 5).Value = 1081
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 19

# Row 80: Republica de Yibuti
$ws.Cells.Item(80, 1).Value = 'Republica de Yibuti'
$ws.Cells.Item(80, 2).Value = 4449
$ws.Cells.Item(80, 3).Value = 8
$ws.Cells.Item(80, 4).Value = 2823
$ws.Cells.Item(80, 5).Value = 1585
$ws.Cells.Item(80, 6).Value = 0
$ws.Cells.Item(80, 7).Value = 3
$ws.Cells.Item(80, 8).Value = 41

# Row 85: Republica de Macedonia
$ws.Cells.Item(85, 1).Value = 'Republica de Macedonia'
$ws.Cells.Item(85, 2).Value = 3895
$ws.Cells.Item(85, 3).Value = 194
$ws.Cells.Item(85, 4).Value = 1705
$ws.Cells.Item(85, 5).Value = 2011
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 8
$ws.Cells.Item(85, 8).Value = 179

# Row 92: Grecia
$ws.Cells.Item(92, 1).Value = 'Grecia'
$ws.Cells.Item(92, 2).Value = 3112
$ws.Cells.Item(92, 3).Value = 4
$ws.Cells.Item(92, 4).Value = 1374
$ws.Cells.Item(92, 5).Value = 1555
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 183

# Row 104: Islandia
$ws.Cells.Item(104, 1).Value = 'Islandia'
$ws.Cells.Item(104, 2).Value = 1808
$ws.Cells.Item(104, 3).Value = 1
$ws.Cells.Item(104, 4).Value = 1794
$ws.Cells.Item(104, 5).Value = 4
$ws.Cells.Item(104, 6).Value = 0
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 10

# Row 121: Sierra Leona
$ws.Cells.Item(121, 1).Value = 'Sierra Leona'
$ws.Cells.Item(121, 2).Value = 1132
$ws.Cells.Item(121, 3).Value = 29
$ws.Cells.Item(121, 4).Value = 670
$ws.Cells.Item(121, 5).Value = 411
$ws.Cells.Item(121, 6).Value = 0
$ws.Cells.Item(121, 7).Value = 0
$ws.Cells.Item(121, 8).Value = 51

# Row 122: Hong Kong
$ws.Cells.Item(122, 1).Value = 'Hong Kong'
$ws.Cells.Item(122, 2).Value = 1110
$ws.Cells.Item(122, 3).Value = 1
$ws.Cells.Item(122, 4).Value = 1061
$ws.Cells.Item(122, 5).Value = 45
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 4

# Row 124: Tunez
$ws.Cells.Item(124, 1).Value = 'Tunez'
$ws.Cells.Item(124, 2).Value = 1094
$ws.Cells.Item(124, 3).Value = 1
$ws.Cells.Item(124, 4).Value = 995
$ws.Cells.Item(124, 5).Value = 50
$ws.Cells.Item(124, 6).Value = 0
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 49

# Row 151: Liberia
$ws.Cells.Item(151, 1).Value = 'Liberia'
$ws.Cells.Item(151, 2).Value = 446
$ws.Cells.Item(151, 3).Value = 25
$ws.Cells.Item(151, 4).Value = 214
$ws.Cells.Item(151, 5).Value = 200
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 32

# Row 152: Taiwan
$ws.Cells.Item(152, 1).Value = 'Taiwan'
$ws.Cells.Item(152, 2).Value = 443
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(152, 4).Value = 431
$ws.Cells.Item(152, 5).Value = 5
$ws.Cells.Item(152, 6).Value = 0
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 7

# Row 160: Birmania
$ws.Cells.Item(160, 1).Value = 'Birmania'
$ws.Cells.Item(160, 2).Value = 261
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(160, 4).Value = 167
$ws.Cells.Item(160, 5).Value = 88
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 0
$ws.Cells.Item(160, 8).Value = 6

# Row 206: Islas Malvinas
$ws.Cells.Item(206, 1).Value = 'Islas Malvinas'
$ws.Cells.Item(206, 2).Value = 13
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 13
$ws.Cells.Item(206, 5).Value = 0
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 0

# Row 207: Groenlandia
$ws.Cells.Item(207, 1).Value = 'Groenlandia'
$ws.Cells.Item(207, 2).Value = 13
$ws.Cells.Item(207, 3).Value = 0
$ws.Cells.Item(207, 4).Value = 13
$ws.Cells.Item(207, 5).Value = 0
$ws.Cells.Item(207, 6).Value = 0
$ws.Cells.Item(207, 7).Value = 0
$ws.Cells.Item(207, 8).Value = 0

# Row 208: Islas Turcas y Caicos
$ws.Cells.Item(208, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(208, 2).Value = 12
$ws.Cells.Item(208, 3).Value = 0
$ws.Cells.Item(208, 4).Value = 11
$ws.Cells.Item(208, 5).Value = 0
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 1

# Row 209: Santa Sede
$ws.Cells.Item(209, 1).Value = 'Santa Sede'
$ws.Cells.Item(209, 2).Value = 12
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 12
$ws.Cells.Item(209, 5).Value = 0
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0

# Row 210: Seychelles
$ws.Cells.Item(210, 1).Value = 'Seychelles'
$ws.Cells.Item(210, 2).Value = 11
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0

# Row 211: Montserrat
$ws.Cells.Item(211, 1).Value = 'Montserrat'
$ws.Cells.Item(211, 2).Value = 11
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 5).Value = 0
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 1

# Row 213: Papua Nueva Guinea
$ws.Cells.Item(213, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(213, 2).Value = 8
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 8
$ws.Cells.Item(213, 5).Value = 0
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 0

# Row 214: Islas Virgenes Britanicas
$ws.Cells.Item(214, 1).Value = 'Islas Virgenes Britanicas'
$ws.Cells.Item(214, 2).Value = 8
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 1
